# Append 6 new product rows (barcode / price pairs) to the bottom of the
# sheet's data table, extending the used range from A2:B101 to A2:B107.
#
# Column A holds long, purely-numeric barcode/EAN identifiers. Those must be
# forced to Text so Excel stores them verbatim (no silent numeric coercion /
# precision loss) - matching how every other cell in column A is already
# stored. Column B's values use a comma decimal separator ("22,4") which
# Excel already treats as plain text on its own, so no extra formatting is
# required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("6971064236152", "22,4"),
    @("6971064236657", "19,4"),
    @("6971064230495", "24,4"),
    @("6971064230709", "79,4"),
    @("3000000016268", "1,08"),
    @("3000000041543", "1,85")
)

$startRow = 102
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $colA = $ws.Cells.Item($r, 1)
    $colB = $ws.Cells.Item($r, 2)

    $colA.NumberFormat = "@"
    $colA.Value = $newRows[$i][0]
    $colB.Value = $newRows[$i][1]
}
